$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# I24 was previously a projected/formula value
# (=I23*(1+AVERAGE(M22:M23))); it is now an actual reported number for
# Mar 29, so replace the formula with the hard-coded actual and restyle
# the cell to match the other "actuals" cells (e.g. I22/I23) instead of
# the "projection" style.
$ws.Range("I23").Copy()
$ws.Range("I24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("I24").Value = 142460

# Move the active selection to reflect the author's last cursor position.
$ws.Range("I25").Select()
